$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4,1).Value = 'ask.fm'
$ws.Cells.Item(4,2).Value = 'https://ask.fm/gorbash'
$ws.Cells.Item(5,1).Value = 'Bandlab'
$ws.Cells.Item(5,2).Value = 'https://www.bandlab.com/api/v1.3/users/gorbash'
$ws.Cells.Item(8,1).Value = 'chaturbate'
$ws.Cells.Item(8,2).Value = 'https://chaturbate.com/gorbash/'
$ws.Cells.Item(10,1).Value = 'Chyoa'
$ws.Cells.Item(10,2).Value = 'https://chyoa.com/user/gorbash'
$ws.Cells.Item(11,1).Value = 'Chomikuj.pl'
$ws.Cells.Item(11,2).Value = 'https://chomikuj.pl/gorbash/'
$ws.Cells.Item(13,1).Value = 'cda.pl'
$ws.Cells.Item(13,2).Value = 'https://www.cda.pl/gorbash'
$ws.Cells.Item(14,1).Value = 'cHEEZburger'
$ws.Cells.Item(14,2).Value = 'https://profile.cheezburger.com/gorbash'
$ws.Cells.Item(15,1).Value = 'cfx.re'
$ws.Cells.Item(15,2).Value = 'https://forum.cfx.re/u/gorbash.json'
$ws.Cells.Item(16,1).Value = 'Codewars'
$ws.Cells.Item(16,2).Value = 'https://www.codewars.com/users/gorbash'
$ws.Cells.Item(19,1).Value = 'Cults3D'
$ws.Cells.Item(19,2).Value = 'https://cults3d.com/en/users/gorbash/creations'
$ws.Cells.Item(20,1).Value = 'diigo'
$ws.Cells.Item(20,2).Value = 'https://www.diigo.com/interact_api/load_profile_info?name=gorbash'
$ws.Cells.Item(22,1).Value = 'Disqus'
$ws.Cells.Item(22,2).Value = 'https://disqus.com/by/gorbash/'
$ws.Cells.Item(23,1).Value = 'Discogs'
$ws.Cells.Item(23,2).Value = 'https://api.discogs.com/users/gorbash'
$ws.Cells.Item(24,1).Value = 'DockerHub'
$ws.Cells.Item(24,2).Value = 'https://hub.docker.com/v2/users/gorbash/'
$ws.Cells.Item(25,1).Value = 'DeviantArt'
$ws.Cells.Item(25,2).Value = 'https://www.deviantart.com/gorbash'
$ws.Cells.Item(28,1).Value = 'Etsy'
$ws.Cells.Item(28,2).Value = 'https://www.etsy.com/people/gorbash'
$ws.Cells.Item(29,1).Value = 'Demotywatory'
$ws.Cells.Item(29,2).Value = 'https://demotywatory.pl/user/gorbash'
$ws.Cells.Item(30,1).Value = 'FatSecret'
$ws.Cells.Item(30,2).Value = 'https://www.fatsecret.com/member/gorbash'
$ws.Cells.Item(33,1).Value = 'Fabswingers'
$ws.Cells.Item(33,2).Value = 'https://www.fabswingers.com/profile/gorbash'
$ws.Cells.Item(35,1).Value = 'Foursquare'
$ws.Cells.Item(35,2).Value = 'https://foursquare.com/gorbash'
$ws.Cells.Item(36,1).Value = 'Filmweb'
$ws.Cells.Item(36,2).Value = 'https://www.filmweb.pl/user/gorbash'
$ws.Cells.Item(38,1).Value = 'Geocaching'
$ws.Cells.Item(38,2).Value = 'https://www.geocaching.com/p/?u=gorbash'
$ws.Cells.Item(39,1).Value = 'Gravatar'
$ws.Cells.Item(39,2).Value = 'https://en.gravatar.com/gorbash.json'
$ws.Cells.Item(40,1).Value = 'GitHub'
$ws.Cells.Item(40,2).Value = 'https://github.com/gorbash'
$ws.Cells.Item(41,1).Value = 'Flipboard'
$ws.Cells.Item(41,2).Value = 'https://flipboard.com/@gorbash'
$ws.Cells.Item(42,1).Value = 'GitLab'
$ws.Cells.Item(42,2).Value = 'https://gitlab.com/gorbash'
$ws.Cells.Item(43,1).Value = 'giters'
$ws.Cells.Item(43,2).Value = 'https://giters.com/gorbash'
$ws.Cells.Item(45,1).Value = 'HudsonRock'
$ws.Cells.Item(45,2).Value = 'https://cavalier.hudsonrock.com/api/json/v2/osint-tools/search-by-username?username=gorbash'
$ws.Cells.Item(46,1).Value = 'IFTTT'
$ws.Cells.Item(46,2).Value = 'https://ifttt.com/p/gorbash'
$ws.Cells.Item(47,1).Value = 'HackerOne'
$ws.Cells.Item(47,2).Value = 'https://hackerone.com/gorbash'
$ws.Cells.Item(50,1).Value = 'inaturalist'
$ws.Cells.Item(50,2).Value = 'https://inaturalist.nz/people/gorbash'
$ws.Cells.Item(51,1).Value = 'issuu'
$ws.Cells.Item(51,2).Value = 'https://issuu.com/gorbash'
$ws.Cells.Item(52,1).Value = 'Instagram2'
$ws.Cells.Item(52,2).Value = 'https://dumpoir.com/v/gorbash'
$ws.Cells.Item(53,1).Value = 'kaggle'
$ws.Cells.Item(53,2).Value = 'https://www.kaggle.com/gorbash'
$ws.Cells.Item(54,1).Value = 'Keybase'
$ws.Cells.Item(54,2).Value = 'https://keybase.io/gorbash'
$ws.Cells.Item(55,1).Value = 'Jeuxvideo'
$ws.Cells.Item(55,2).Value = 'https://www.jeuxvideo.com/profil/gorbash?mode=infos'
$ws.Cells.Item(56,1).Value = 'Internet Archive User Search'
$ws.Cells.Item(56,2).Value = 'https://archive.org/advancedsearch.php?q=gorbash&output=json'
$ws.Cells.Item(57,1).Value = 'Kongregate'
$ws.Cells.Item(57,2).Value = 'https://www.kongregate.com/accounts/gorbash'
$ws.Cells.Item(59,1).Value = 'InkBunny'
$ws.Cells.Item(59,2).Value = 'https://inkbunny.net/gorbash'
$ws.Cells.Item(64,1).Value = 'MCUUID (Minecraft)'
$ws.Cells.Item(64,2).Value = 'https://playerdb.co/api/player/minecraft/gorbash'
$ws.Cells.Item(65,1).Value = 'Livejournal'
$ws.Cells.Item(65,2).Value = 'https://gorbash.livejournal.com'
$ws.Cells.Item(66,1).Value = 'Mastodon API'
$ws.Cells.Item(66,2).Value = 'https://mastodon.social/api/v2/search?q=gorbash&limit=1&type=accounts'
$ws.Cells.Item(68,1).Value = 'MySpace'
$ws.Cells.Item(68,2).Value = 'https://myspace.com/gorbash'
$ws.Cells.Item(69,1).Value = 'MyAnimeList'
$ws.Cells.Item(69,2).Value = 'https://myanimelist.net/profile/gorbash'
$ws.Cells.Item(71,1).Value = 'palnet'
$ws.Cells.Item(71,2).Value = 'https://www.palnet.io/@gorbash/'
$ws.Cells.Item(73,1).Value = 'pikabu'
$ws.Cells.Item(73,2).Value = 'https://pikabu.ru/@gorbash'
$ws.Cells.Item(74,1).Value = 'Pinterest'
$ws.Cells.Item(74,2).Value = 'https://www.pinterest.com/gorbash/'
$ws.Cells.Item(75,1).Value = 'Pokerstrategy'
$ws.Cells.Item(75,2).Value = 'http://www.pokerstrategy.net/user/gorbash/profile/'
$ws.Cells.Item(76,1).Value = 'Periscope'
$ws.Cells.Item(76,2).Value = 'https://www.periscope.tv/gorbash'
$ws.Cells.Item(77,1).Value = 'prv.pl'
$ws.Cells.Item(77,2).Value = 'https://www.prv.pl/osoba/gorbash'
$ws.Cells.Item(79,1).Value = 'public'
$ws.Cells.Item(79,2).Value = 'https://public.com/@gorbash'
$ws.Cells.Item(80,1).Value = 'Roblox'
$ws.Cells.Item(80,2).Value = 'https://auth.roblox.com/v1/usernames/validate?username=gorbash&birthday=2019-12-31T23:00:00.000Z'
$ws.Cells.Item(82,1).Value = 'MCName (Minecraft)'
$ws.Cells.Item(82,2).Value = 'https://mcname.info/en/search?q=gorbash'
$ws.Cells.Item(83,1).Value = 'RumbleUser'
$ws.Cells.Item(83,2).Value = 'https://rumble.com/user/gorbash'
$ws.Cells.Item(84,1).Value = 'slideshare'
$ws.Cells.Item(84,2).Value = 'https://www.slideshare.net/gorbash'
$ws.Cells.Item(86,1).Value = 'SFD'
$ws.Cells.Item(86,2).Value = 'https://www.sfd.pl/profile/gorbash'
$ws.Cells.Item(88,1).Value = 'sofurry'
$ws.Cells.Item(88,2).Value = 'https://gorbash.sofurry.com'
$ws.Cells.Item(89,1).Value = 'SoundCloud'
$ws.Cells.Item(89,2).Value = 'https://soundcloud.com/gorbash'
$ws.Cells.Item(90,1).Value = 'Sourceforge'
$ws.Cells.Item(90,2).Value = 'https://sourceforge.net/u/gorbash/profile'
$ws.Cells.Item(91,1).Value = 'Spotify'
$ws.Cells.Item(91,2).Value = 'https://open.spotify.com/user/gorbash'
$ws.Cells.Item(92,1).Value = 'Steam'
$ws.Cells.Item(92,2).Value = 'https://steamcommunity.com/id/gorbash'
$ws.Cells.Item(93,1).Value = 'themeforest'
$ws.Cells.Item(93,2).Value = 'https://themeforest.net/user/gorbash'
$ws.Cells.Item(95,1).Value = 'Telegram'
$ws.Cells.Item(95,2).Value = 'https://t.me/gorbash'
$ws.Cells.Item(96,1).Value = 'TikTok'
$ws.Cells.Item(96,2).Value = 'https://www.tiktok.com/oembed?url=https://www.tiktok.com/@gorbash'
$ws.Cells.Item(97,1).Value = 'Trello'
$ws.Cells.Item(97,2).Value = 'https://trello.com/1/Members/gorbash?fields=activityBlocked%2CavatarUrl%2Cbio%2CbioData%2Cconfirmed%2CfullName%2CidEnterprise%2CidMemberReferrer%2Cinitials%2CmemberType%2CnonPublic%2Cproducts%2Curl%2Cusername'
$ws.Cells.Item(98,1).Value = 'Twitter archived profile'
$ws.Cells.Item(98,2).Value = 'http://archive.org/wayback/available?url=https://twitter.com/gorbash'
$ws.Cells.Item(101,1).Value = 'Twitch'
$ws.Cells.Item(101,2).Value = 'https://twitchtracker.com/gorbash'
$ws.Cells.Item(102,1).Value = 'Twitter'
$ws.Cells.Item(102,2).Value = 'https://nitter.privacydev.net/gorbash'
$ws.Cells.Item(103,1).Value = 'Venmo'
$ws.Cells.Item(103,2).Value = 'https://account.venmo.com/u/gorbash'
$ws.Cells.Item(105,1).Value = 'untappd'
$ws.Cells.Item(105,2).Value = 'https://untappd.com/user/gorbash/'
$ws.Cells.Item(106,1).Value = 'wattpad'
$ws.Cells.Item(106,2).Value = 'https://www.wattpad.com/api/v3/users/gorbash?fields=username%2Cname%2Cdescription%2Cavatar%2CbackgroundUrl%2CcreateDate%2Clocation%2Cfollowing%2CfollowingRequest%2CnumFollowing%2Cfollower%2CfollowerRequest%2CnumFollowers%2CnumLists%2CnumStoriesPublished%2CvotesReceived%2Cfacebook%2Ctwitter%2Cwebsite%2Csmashwords%2Chighlight_colour%2Chtml_enabled%2Cverified%2Cambassador%2Cwattpad_squad%2Cis_staff%2Cprograms(wattpad_stars)%2CisPrivate%2CisMuted%2CexternalId%2Cnotes'
$ws.Cells.Item(109,1).Value = 'WordPress Support'
$ws.Cells.Item(109,2).Value = 'https://wordpress.org/support/users/gorbash/'
$ws.Cells.Item(110,1).Value = 'VK'
$ws.Cells.Item(110,2).Value = 'https://vk.com/gorbash'
$ws.Cells.Item(111,1).Value = 'YouTube User2'
$ws.Cells.Item(111,2).Value = 'https://www.youtube.com/@gorbash'
$ws.Cells.Item(113,1).Value = 'xHamster'
$ws.Cells.Item(113,2).Value = 'https://xhamster.com/users/gorbash'

$ws.Range("A118:B118").Delete()

Write-Output "done"